$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Fix the title of the existing "EFA-final" sheet: drop " (Final)"
# ------------------------------------------------------------------
$final = $wb.Worksheets.Item("EFA-final")
$final.Range("A1").Value = "Exploratory Factor Analysis EFA with n=3 factors"

# ------------------------------------------------------------------
# 2) Create the new "EFA-final-alt" sheet as a copy of "EFA-final" so
#    it inherits the same layout / number formats / column widths,
#    then overwrite its contents with the alternative-model results.
# ------------------------------------------------------------------
$final.Copy($null, $final)
$newSheet = $wb.Worksheets.Item("EFA-final (2)")
$newSheet.Name = "EFA-final-alt"

$ws = $newSheet

# Remove the old merged footer cell and clear every cell so no
# leftover data from the longer (17-item) "EFA-final" table remains.
$ws.Range("A63:J67").UnMerge()
$ws.Range("A1:J68").ClearContents()

# ------------------------------------------------------------------
# 3) Write the alternative EFA (n=3 factors, 13 items) results.
# ------------------------------------------------------------------
    $ws.Range("A1").Value = "Exploratory Factor Analysis Alternative EFA with n=3 factors"
    $ws.Range("A2").Value = "  "
    $ws.Range("A3").Value = "  "
    $ws.Range("A4").Value = "Standardized loadings"
    $ws.Range("B5").Value = "ML2"
    $ws.Range("C5").Value = "ML1"
    $ws.Range("D5").Value = "ML3"
    $ws.Range("A6").Value = "Item3"
    $ws.Range("B6").Value = -0.038
    $ws.Range("C6").Value = 0.8356
    $ws.Range("D6").Value = -0.0612
    $ws.Range("A7").Value = "Item4"
    $ws.Range("B7").Value = -0.1972
    $ws.Range("C7").Value = 0.8973
    $ws.Range("D7").Value = 0.1486
    $ws.Range("A8").Value = "Item5"
    $ws.Range("B8").Value = 0.0925
    $ws.Range("C8").Value = 0.7538
    $ws.Range("D8").Value = -0.0344
    $ws.Range("A9").Value = "Item6"
    $ws.Range("B9").Value = 0.314
    $ws.Range("C9").Value = 0.5791
    $ws.Range("D9").Value = -0.0743
    $ws.Range("A10").Value = "Item12"
    $ws.Range("B10").Value = 0.4548
    $ws.Range("C10").Value = -0.0238
    $ws.Range("D10").Value = 0.0538
    $ws.Range("A11").Value = "Item13"
    $ws.Range("B11").Value = 0.9056
    $ws.Range("C11").Value = 0.0387
    $ws.Range("D11").Value = -0.1693
    $ws.Range("A12").Value = "Item14"
    $ws.Range("B12").Value = 0.5068
    $ws.Range("C12").Value = -0.0117
    $ws.Range("D12").Value = 0.0761
    $ws.Range("A13").Value = "Item15"
    $ws.Range("B13").Value = 0.8842
    $ws.Range("C13").Value = -0.0926
    $ws.Range("D13").Value = -0.069
    $ws.Range("A14").Value = "Item17"
    $ws.Range("B14").Value = 0.6217
    $ws.Range("C14").Value = 0.0048
    $ws.Range("D14").Value = 0.1145
    $ws.Range("A15").Value = "Item18"
    $ws.Range("B15").Value = 0.5713
    $ws.Range("C15").Value = 0.1115
    $ws.Range("D15").Value = 0.0924
    $ws.Range("A16").Value = "Item19"
    $ws.Range("B16").Value = 0.026
    $ws.Range("C16").Value = 0.0203
    $ws.Range("D16").Value = 0.7114
    $ws.Range("A17").Value = "Item21"
    $ws.Range("B17").Value = -0.0458
    $ws.Range("C17").Value = -0.0432
    $ws.Range("D17").Value = 0.813
    $ws.Range("A18").Value = "Item23"
    $ws.Range("B18").Value = 0.2154
    $ws.Range("C18").Value = 0.129
    $ws.Range("D18").Value = 0.3282
    $ws.Range("A19").Value = "  "
    $ws.Range("A20").Value = "  "
    $ws.Range("A21").Value = "Factor correlations"
    $ws.Range("B22").Value = "ML2"
    $ws.Range("C22").Value = "ML1"
    $ws.Range("D22").Value = "ML3"
    $ws.Range("A23").Value = "ML2"
    $ws.Range("B23").Value = 1
    $ws.Range("C23").Value = 0.7419
    $ws.Range("D23").Value = 0.7815
    $ws.Range("A24").Value = "ML1"
    $ws.Range("B24").Value = 0.7419
    $ws.Range("C24").Value = 1
    $ws.Range("D24").Value = 0.7347
    $ws.Range("A25").Value = "ML3"
    $ws.Range("B25").Value = 0.7815
    $ws.Range("C25").Value = 0.7347
    $ws.Range("D25").Value = 1
    $ws.Range("A26").Value = "  "
    $ws.Range("A27").Value = "  "
    $ws.Range("A28").Value = "Measures of factor score adequacy"
    $ws.Range("B29").Value = "Vaccounted.ML2"
    $ws.Range("C29").Value = "Vaccounted.ML1"
    $ws.Range("D29").Value = "Vaccounted.ML3"
    $ws.Range("A30").Value = "SS loadings"
    $ws.Range("B30").Value = 3.0067
    $ws.Range("C30").Value = 2.5279
    $ws.Range("D30").Value = 1.3785
    $ws.Range("A31").Value = "Proportion Var"
    $ws.Range("B31").Value = 0.2313
    $ws.Range("C31").Value = 0.1945
    $ws.Range("D31").Value = 0.106
    $ws.Range("A32").Value = "Cumulative Var"
    $ws.Range("B32").Value = 0.2313
    $ws.Range("C32").Value = 0.4257
    $ws.Range("D32").Value = 0.5318
    $ws.Range("A33").Value = "Proportion Explained"
    $ws.Range("B33").Value = 0.4349
    $ws.Range("C33").Value = 0.3657
    $ws.Range("D33").Value = 0.1994
    $ws.Range("A34").Value = "Cumulative Proportion"
    $ws.Range("B34").Value = 0.4349
    $ws.Range("C34").Value = 0.8006
    $ws.Range("D34").Value = 1
    $ws.Range("A35").Value = "  "
    $ws.Range("A36").Value = "  "
    $ws.Range("A37").Value = "Item complexity"
    $ws.Range("B38").Value = "Value"
    $ws.Range("A39").Value = "Item3"
    $ws.Range("B39").Value = 1.0149
    $ws.Range("A40").Value = "Item4"
    $ws.Range("B40").Value = 1.1536
    $ws.Range("A41").Value = "Item5"
    $ws.Range("B41").Value = 1.0343
    $ws.Range("A42").Value = "Item6"
    $ws.Range("B42").Value = 1.5804
    $ws.Range("A43").Value = "Item12"
    $ws.Range("B43").Value = 1.0335
    $ws.Range("A44").Value = "Item13"
    $ws.Range("B44").Value = 1.0736
    $ws.Range("A45").Value = "Item14"
    $ws.Range("B45").Value = 1.0462
    $ws.Range("A46").Value = "Item15"
    $ws.Range("B46").Value = 1.0343
    $ws.Range("A47").Value = "Item17"
    $ws.Range("B47").Value = 1.0679
    $ws.Range("A48").Value = "Item18"
    $ws.Range("B48").Value = 1.1301
    $ws.Range("A49").Value = "Item19"
    $ws.Range("B49").Value = 1.0043
    $ws.Range("A50").Value = "Item21"
    $ws.Range("B50").Value = 1.012
    $ws.Range("A51").Value = "Item23"
    $ws.Range("B51").Value = 2.078
    $ws.Range("A52").Value = "  "
    $ws.Range("A53").Value = "  "
    $ws.Range("A54").Value = "Extra information"
    $ws.Range("A55").Value = "Mean item complexity = 1.17409070887675`nThe degrees of freedom for the model are = 42`nThe objective function was = 0.426582512916918`nThe Chi Square of the model is = 107.887316779403`nThe root mean square of the residuals (RMSR) is = 0.0426610066000907`nTucker Lewis Index of factoring reliability is =  0.901885178102094`nRMSEA index =  0.0854395900327584`nThe 90% confidence intervals of RMSEA are lower = 0.0716646306244335 and upper = 0.0999634220637349`nBIC = -90.8695948946509`n"
    $ws.Range("A56").Value = "  "
    $ws.Range("A57").Value = "  "
    $ws.Range("A58").Value = "  "
    $ws.Range("A59").Value = "  "
    $ws.Range("A60").Value = "  "

# ------------------------------------------------------------------
# 4) Re-create the merged footer cell at its new (shorter) location.
# ------------------------------------------------------------------
$ws.Range("A55:J59").Merge()
